$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.203.79'
$ws.Range("E2").Value = '  -1.67%  '
$ws.Range("D3").Value = '3.792.97'
$ws.Range("E3").Value = '  +3.70%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '''619.62'
$ws.Range("E5").Value = '  +2.97%  '
$ws.Range("E6").Value = '  -4.29%  '
$ws.Range("D7").Value = '3.788.09'
$ws.Range("D9").Value = '''0.537'
$ws.Range("E9").Value = '  +0.04%  '
$ws.Range("D10").Value = '''0.171'
$ws.Range("E10").Value = '  +3.28%  '
$ws.Range("E11").Value = '  -3.56%  '
$ws.Range("E12").Value = '  -1.35%  '
$ws.Range("D13").Value = '''41.18'
$ws.Range("E13").Value = '  +2.58%  '
$ws.Range("D14").Value = '''0.0000260'
$ws.Range("E14").Value = '  +0.36%  '
$ws.Range("D15").Value = '4.426.07'
$ws.Range("E15").Value = '  +3.60%  '
$ws.Range("D16").Value = '3.792.28'
$ws.Range("E16").Value = '  +3.45%  '
$ws.Range("D17").Value = '70.261.12'
$ws.Range("E17").Value = '  -1.57%  '
$ws.Range("D19").Value = '''7.62'
$ws.Range("E19").Value = '  +0.72%  '
$ws.Range("D20").Value = '''515.74'
$ws.Range("E20").Value = '  +0.18%  '
$ws.Range("D21").Value = '''16.77'
$ws.Range("E21").Value = '  -3.13%  '
$ws.Range("D22").Value = '''9.64'
$ws.Range("E22").Value = '  +3.28%  '
$ws.Range("E23").Value = '  -2.70%  '
$ws.Range("E24").Value = '  +4.81%  '
$ws.Range("D25").Value = '''88.24'
$ws.Range("E25").Value = '  -0.39%  '
$ws.Range("D26").Value = '''13.31'
$ws.Range("E26").Value = '  -1.93%  '
$ws.Range("D27").Value = '''11.22'
$ws.Range("E27").Value = '  +3.24%  '
$ws.Range("E28").Value = '  +22.77%  '
$ws.Range("D29").Value = '''0.998'
$ws.Range("E29").Value = '  -0.11%  '
$ws.Range("E30").Value = '  -2.35%  '
$ws.Range("E31").Value = '  +2.84%  '
$ws.Range("D32").Value = '''7.83'
$ws.Range("E32").Value = '  -5.28%  '
$ws.Range("D33").Value = '''31.76'
$ws.Range("E34").Value = '  -1.94%  '
$ws.Range("D35").Value = '''1.00'
$ws.Range("E35").Value = '  +0.04%  '
$ws.Range("E36").Value = '  +1.18%  '
$ws.Range("E37").Value = '  +2.28%  '
$ws.Range("D38").Value = '''0.340'
$ws.Range("E38").Value = '  +0.61%  '
$ws.Range("E39").Value = '  +1.85%  '
$ws.Range("E40").Value = '  +3.41%  '
$ws.Range("D41").Value = '''51.13'
$ws.Range("E41").Value = '  +0.25%  '
$ws.Range("D42").Value = '''44.45'
$ws.Range("E42").Value = '  -6.07%  '
$ws.Range("E43").Value = '  -1.45%  '
$ws.Range("D44").Value = '''423.32'
$ws.Range("E44").Value = '  +3.97%  '
$ws.Range("D45").Value = '3.066.36'
$ws.Range("E45").Value = '  -2.96%  '
$ws.Range("E46").Value = '  -2.44%  '
$ws.Range("E47").Value = '  -0.33%  '
$ws.Range("D48").Value = '''27.65'
$ws.Range("E48").Value = '  -1.91%  '
$ws.Range("D49").Value = '''136.14'
$ws.Range("E49").Value = '  +0.94%  '
$ws.Range("E51").Value = '  -0.30%  '
